# Insert a new weekly record as row 83 on the single data sheet,
# pushing the existing rows 83-140 down to 84-141 (dimension grows to A1:R141).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 83 (shifts rows 83..140 down to 84..141).
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new record's data.
$ws.Cells.Item(83, 1).Value = 5
$ws.Cells.Item(83, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(83, 3).Value = "Maule"
$ws.Cells.Item(83, 4).Value = 44582
$ws.Cells.Item(83, 5).Value = 7
$ws.Cells.Item(83, 6).Value = 100112031
$ws.Cells.Item(83, 7).Value = "Poroto verde"
$ws.Cells.Item(83, 8).Value = "Sin especificar"
$ws.Cells.Item(83, 9).Value = "Primera"
$ws.Cells.Item(83, 10).Value = 60
$ws.Cells.Item(83, 11).Value = 30000
$ws.Cells.Item(83, 12).Value = 30000
$ws.Cells.Item(83, 13).Value = 30000
$ws.Cells.Item(83, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(83, 15).Value = "Región del Maule"
$ws.Cells.Item(83, 16).Value = 1200
$ws.Cells.Item(83, 17).Value = 25
$ws.Cells.Item(83, 18).Value = "Hortaliza"
